# Update "想去人数" (want-to-go count) figures across sheets, reflecting
# the latest data pull (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1421
$ws1.Range("F4").Value = 1058
$ws1.Range("F8").Value = 214
$ws1.Range("F10").Value = 57
$ws1.Range("F12").Value = 132
$ws1.Range("F13").Value = 1699
$ws1.Range("F14").Value = 339
$ws1.Range("F18").Value = 402
$ws1.Range("F20").Value = 2
$ws1.Range("F21").Value = 643
$ws1.Range("F23").Value = 227
$ws1.Range("F24").Value = 944
$ws1.Range("F25").Value = 51
$ws1.Range("F26").Value = 1497
$ws1.Range("F27").Value = 225

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 119
$ws2.Range("F3").Value = 33
$ws2.Range("F9").Value = 38

# --- 本地生活 sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 390

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1421
$ws4.Range("F5").Value = 1058
$ws4.Range("F6").Value = 119
$ws4.Range("F7").Value = 33
$ws4.Range("F12").Value = 214
$ws4.Range("F14").Value = 57
$ws4.Range("F16").Value = 132
$ws4.Range("F17").Value = 1699
$ws4.Range("F19").Value = 339
$ws4.Range("F23").Value = 402
$ws4.Range("F26").Value = 2
$ws4.Range("F29").Value = 643
$ws4.Range("F30").Value = 38
$ws4.Range("F35").Value = 227
$ws4.Range("F36").Value = 944
$ws4.Range("F37").Value = 51
$ws4.Range("F38").Value = 1497
$ws4.Range("F39").Value = 225
